$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 8 (year 2025) metrics with the latest figures
$ws.Range("C8").Value = 992
$ws.Range("D8").Value = 164
$ws.Range("E8").Value = 828
$ws.Range("F8").Value = 6.726825266611977
$ws.Range("G8").Value = 83.46774193548387
$ws.Range("H8").Value = 16.53225806451613
